$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 333.0
$ws.Range("B3").Value = 255.66
$ws.Range("B4").Value = 77.34
$ws.Range("B5").Value = 32.75
